$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-03 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-04 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("92÷6=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=8, 4", 2) | Out-Null
$d.Content.Find.Execute("83÷6=13, 5", $true, $false, $false, $false, $false, $true, 1, $false, "34÷2=17, 0", 2) | Out-Null
$d.Content.Find.Execute("40÷4=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "54÷3=18, 0", 2) | Out-Null
$d.Content.Find.Execute("76÷6=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "53÷8=6, 5", 2) | Out-Null
$d.Content.Find.Execute("58÷7=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=9, 3", 2) | Out-Null
$d.Content.Find.Execute("22÷4=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "78÷9=8, 6", 2) | Out-Null
$d.Content.Find.Execute("77÷6=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "85÷6=14, 1", 2) | Out-Null
$d.Content.Find.Execute("85÷4=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "42÷7=6, 0", 2) | Out-Null
$d.Content.Find.Execute("99÷2=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$d.Content.Find.Execute("90÷6=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=9, 3", 2) | Out-Null
$d.Content.Find.Execute("48÷4=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷3=19, 2", 2) | Out-Null
$d.Content.Find.Execute("97÷2=48, 1", $true, $false, $false, $false, $false, $true, 1, $false, "12÷4=3, 0", 2) | Out-Null
$d.Content.Find.Execute("74÷2=37, 0", $true, $false, $false, $false, $false, $true, 1, $false, "30÷6=5, 0", 2) | Out-Null
$d.Content.Find.Execute("75÷6=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "50÷4=12, 2", 2) | Out-Null
$d.Content.Find.Execute("20÷3=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "24÷8=3, 0", 2) | Out-Null
$d.Content.Find.Execute("97÷6=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "65÷3=21, 2", 2) | Out-Null
$d.Content.Find.Execute("22÷3=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "64÷9=7, 1", 2) | Out-Null
$d.Content.Find.Execute("35÷5=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "10÷9=1, 1", 2) | Out-Null
$d.Content.Find.Execute("98÷5=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "77÷5=15, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷5=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "40÷4=10, 0", 2) | Out-Null
$d.Content.Find.Execute("82÷3=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷5=4, 1", 2) | Out-Null
$d.Content.Find.Execute("23÷9=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=7, 6", 2) | Out-Null
$d.Content.Find.Execute("39÷2=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "95÷3=31, 2", 2) | Out-Null
$d.Content.Find.Execute("30÷7=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷5=3, 3", 2) | Out-Null
$d.Content.Find.Execute("91÷6=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "52÷8=6, 4", 2) | Out-Null
